# Insert a new data row at row 429 (pushes existing rows 429-528 down to 430-529)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(429).Insert()

# Populate the newly inserted row with its values
$ws.Range("A429").Value = 7
$ws.Range("B429").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C429").Value = "Ñuble"
$ws.Range("D429").Value = 45204
$ws.Range("D429").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E429").Value = 16
$ws.Range("F429").Value = 100112003
$ws.Range("G429").Value = "Ajo"
$ws.Range("H429").Value = "Chino"
$ws.Range("I429").Value = "Primera"
$ws.Range("J429").Value = 80
$ws.Range("K429").Value = 20000
$ws.Range("L429").Value = 20000
$ws.Range("M429").Value = 20000
$ws.Range("N429").Value = "$/caja 10 kilos"
$ws.Range("O429").Value = "China"
$ws.Range("P429").Value = 2000
$ws.Range("Q429").Value = 10
$ws.Range("R429").Value = "Hortaliza"
